$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.986.51"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "2.242.20"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "98.61"
$ws.Range("E5").Value = "  +17.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "271.18"
$ws.Range("E6").Value = "  +4.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.638"
$ws.Range("E9").Value = "  +6.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.91"
$ws.Range("E10").Value = "  +6.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0942"
$ws.Range("E11").Value = "  +2.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.38"
$ws.Range("E12").Value = "  +16.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("E14").Value = "  +6.89%  "
$ws.Range("D15").Value = "2.576.88"
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("E16").Value = "  +4.98%  "
$ws.Range("D17").Value = "2.242.25"
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").Value = "43.954.08"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.21"
$ws.Range("E20").Value = "  +4.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.01"
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.42"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.18"
$ws.Range("E24").Value = "  +2.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.50"
$ws.Range("E26").Value = "  +7.98%  "
$ws.Range("E27").Value = "  +11.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.54"
$ws.Range("E28").Value = "  +2.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.91"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.72"
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0921"
$ws.Range("E32").Value = "  +7.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.08"
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.59"
$ws.Range("E34").Value = "  +4.63%  "
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.112"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("E38").Value = "  -4.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.69"
$ws.Range("E39").Value = "  +30.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.251"
$ws.Range("E40").Value = "  +25.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.62"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("E42").Value = "  +4.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.18"
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.43"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("E45").Value = "  +4.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.49"
$ws.Range("E46").Value = "  +2.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.42"
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("E48").Value = "  +4.08%  "
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.437"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").Value = "2.460.49"
$ws.Range("E51").Value = "  +2.08%  "
